# Weekly fruit/veggie price update: a new record is inserted as the new
# row 51 (for 2021-10-19), pushing all subsequent rows (old 51..77) down
# by one (new 52..78).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51; everything currently at 51..77 shifts
# down to 52..78 (carrying its data/format along, matching the diff).
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new weekly record.
$ws.Range("A51").Value2 = 11
$ws.Range("B51").Value2 = "Vega Monumental Concepción"
$ws.Range("C51").Value2 = "Bíobío"
$ws.Range("D51").Value2 = 44488
$ws.Range("E51").Value2 = 8
$ws.Range("F51").Value2 = 100112043
$ws.Range("G51").Value2 = "Pepino ensalada"
$ws.Range("H51").Value2 = "Sin especificar"
$ws.Range("I51").Value2 = "Primera"
$ws.Range("J51").Value2 = 100
$ws.Range("K51").Value2 = 7000
$ws.Range("L51").Value2 = 7500
$ws.Range("M51").Value2 = 7250
$ws.Range("N51").Value2 = "`$/caja 60 unidades"
$ws.Range("O51").Value2 = "Región de Arica y Parinacota"
$ws.Range("P51").Value2 = 121
$ws.Range("Q51").Value2 = 60
$ws.Range("R51").Value2 = "Hortaliza"
